$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Register")

# Update the password test value in A18 ("Admin@12345" -> "Admin1@567")
$ws.Range("A18").Value = "Admin1@567"

# Bring the Register sheet to front and restore its saved view state:
# scrolled so row 7 is the top-left visible row, with B21 as the
# selected cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B21").Select()
